$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1483.0435
$ws.Range("I98").Value = 1483.0435
$ws.Range("K98").Value = 1483.0435
$ws.Range("M98").Value = 14.95650000000001
$ws.Range("H122").Value = 1483.0435
$ws.Range("I122").Value = 1483.0435
$ws.Range("K122").Value = 4449.1305
$ws.Range("M122").Value = -1999.1305
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3013.0625
$ws.Range("I63").Value = 1727.2727
$ws.Range("J63").Value = 5841.8
$ws.Range("K63").Value = 1727.2727
$ws.Range("L63").Value = 5841.8
$ws.Range("M63").Value = -1041.2727
$ws.Range("N63").Value = -7213.8
$ws.Range("H66").Value = 3013.0625
$ws.Range("I66").Value = 1727.2727
$ws.Range("J66").Value = 5841.8
$ws.Range("K66").Value = 8636.363499999999
$ws.Range("L66").Value = 29209
$ws.Range("M66").Value = -5204.363499999999
$ws.Range("N66").Value = -36073
$ws.Range("H132").Value = 1696.5454
$ws.Range("I132").Value = 1705.8334
$ws.Range("K132").Value = 5117.5002
$ws.Range("M132").Value = -2587.5002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3518.7742
$ws.Range("I105").Value = 3338.25
$ws.Range("K105").Value = 3338.25
$ws.Range("M105").Value = -1591.25
$ws.Range("H107").Value = 1692.45
$ws.Range("I107").Value = 1463.7222
$ws.Range("J107").Value = 3751
$ws.Range("K107").Value = 1463.7222
$ws.Range("L107").Value = 3751
$ws.Range("M107").Value = 456.2778000000001
$ws.Range("N107").Value = -7591
$ws.Range("H132").Value = 83296.914
$ws.Range("J132").Value = 83296.914
$ws.Range("L132").Value = 83296.914
$ws.Range("N132").Value = -93416.914
$ws.Range("H134").Value = 13147.233
$ws.Range("I134").Value = 7390.905
$ws.Range("K134").Value = 22172.715
$ws.Range("M134").Value = -19637.715
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2731.5925
$ws.Range("I31").Value = 1484.8667
$ws.Range("K31").Value = 1484.8667
$ws.Range("M31").Value = -1189.8667
$ws.Range("H34").Value = 2731.5925
$ws.Range("I34").Value = 1484.8667
$ws.Range("K34").Value = 1484.8667
$ws.Range("M34").Value = -1282.8667
$ws.Range("H122").Value = 1953.6957
$ws.Range("I122").Value = 1144.625
$ws.Range("J122").Value = 3803
$ws.Range("K122").Value = 3433.875
$ws.Range("L122").Value = 11409
$ws.Range("M122").Value = -983.875
$ws.Range("N122").Value = -16309
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1600.963
$ws.Range("I5").Value = 715.3570999999999
$ws.Range("J5").Value = 2554.6924
$ws.Range("K5").Value = 2146.0713
$ws.Range("L5").Value = 7664.0772
$ws.Range("M5").Value = -2034.0713
$ws.Range("N5").Value = -7888.0772
$ws.Range("H110").Value = 26949.75
$ws.Range("I110").Value = 26949.75
$ws.Range("K110").Value = 80849.25
$ws.Range("M110").Value = -76759.25
$ws.Range("H129").Value = 2041.8948
$ws.Range("I129").Value = 679.625
$ws.Range("J129").Value = 3032.6365
$ws.Range("K129").Value = 2038.875
$ws.Range("L129").Value = 9097.9095
$ws.Range("M129").Value = 2961.125
$ws.Range("N129").Value = -19097.9095
$ws.Range("H135").Value = 1600.963
$ws.Range("I135").Value = 715.3570999999999
$ws.Range("J135").Value = 2554.6924
$ws.Range("K135").Value = 6438.2139
$ws.Range("L135").Value = 22992.2316
$ws.Range("M135").Value = -3903.2139
$ws.Range("N135").Value = -28062.2316
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 69.8421
$ws.Range("I2").Value = 25.666666
$ws.Range("J2").Value = 145.57143
$ws.Range("K2").Value = 25.666666
$ws.Range("L2").Value = 145.57143
$ws.Range("M2").Value = 87.33333400000001
$ws.Range("N2").Value = -371.57143
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2178.8235
$ws.Range("I82").Value = 1712.5454
$ws.Range("J82").Value = 3033.6667
$ws.Range("K82").Value = 1712.5454
$ws.Range("L82").Value = 3033.6667
$ws.Range("M82").Value = -1351.5454
$ws.Range("N82").Value = -3755.6667
$ws.Range("H85").Value = 2178.8235
$ws.Range("I85").Value = 1712.5454
$ws.Range("J85").Value = 3033.6667
$ws.Range("K85").Value = 1712.5454
$ws.Range("L85").Value = 3033.6667
$ws.Range("M85").Value = -464.5454
$ws.Range("N85").Value = -5529.6667
$ws.Range("H93").Value = 3976.3
$ws.Range("I93").Value = 4769.1055
$ws.Range("K93").Value = 4769.1055
$ws.Range("M93").Value = -3521.1055
$ws.Range("H96").Value = 49998.5
$ws.Range("J96").Value = 49998.5
$ws.Range("L96").Value = 49998.5
$ws.Range("N96").Value = -55490.5
$ws.Range("H132").Value = 5967.2915
$ws.Range("I132").Value = 5706.684
$ws.Range("K132").Value = 17120.052
$ws.Range("M132").Value = -14590.052
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 49995
$ws.Range("I40").Value = 49995
$ws.Range("K40").Value = 49995
$ws.Range("M40").Value = -49846
$ws.Range("H87").Value = 14000
$ws.Range("I87").Value = 14000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 14000
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -12752
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 14000
$ws.Range("I90").Value = 14000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 42000
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -35760
$ws.Range("N90").ClearContents()
$ws.Range("H100").Value = 917.6
$ws.Range("I100").Value = 1047.2354
$ws.Range("K100").Value = 2094.4708
$ws.Range("M100").Value = -1553.4708
$ws.Range("H107").Value = 1446.4736
$ws.Range("I107").Value = 1310.7037
$ws.Range("J107").Value = 1779.7273
$ws.Range("K107").Value = 3932.1111
$ws.Range("L107").Value = 5339.1819
$ws.Range("M107").Value = -2012.1111
$ws.Range("N107").Value = -9179.1819
$ws.Range("H122").Value = 1275.3793
$ws.Range("I122").Value = 507.48
$ws.Range("J122").Value = 6074.75
$ws.Range("K122").Value = 1522.44
$ws.Range("L122").Value = 18224.25
$ws.Range("M122").Value = 927.5599999999999
$ws.Range("N122").Value = -23124.25
$ws.Range("H132").Value = 135460.73
$ws.Range("I132").Value = 180444.81
$ws.Range("K132").Value = 541334.4299999999
$ws.Range("M132").Value = -538804.4299999999
$ws.Range("H136").Value = 3449772.5
$ws.Range("I136").Value = 4762848.5
$ws.Range("K136").Value = 14288545.5
$ws.Range("M136").Value = -14285995.5
